# Applies the cryptos.xlsx price/volume refresh described by the commit:
# "Updated cryptos list on Tue Nov 14 23:47:48 UTC 2023 with GitHub Actions"
#
# All Coin/Link/Price/Volume(1h) cells are plain text in this sheet (Price looks
# numeric but is stored as text, e.g. "35.619.71" / "0.0729"). Two Price cells
# (D7, D27) are new values that end in a trailing zero ("57.10", "163.50") which
# Excel would otherwise silently normalize away if the cell were left in General
# number format, so those two are explicitly forced to Text format first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.619.71"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "1.982.47"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "242.44"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "0.637"
$ws.Range("E6").Value = "  -5.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.10"
$ws.Range("E7").Value = "  +9.40%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "59.62"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "0.359"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "0.0729"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").Value = "0.921"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "14.02"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "2.274.55"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "1.983.21"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").Value = "17.18"
$ws.Range("E18").Value = "  +5.65%  "
$ws.Range("D19").Value = "35.485.86"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "70.65"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "233.42"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  -3.45%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +10.29%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -3.17%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.50"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "9.12"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "19.41"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "4.79"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "0.0589"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "0.0898"
$ws.Range("E34").Value = "  +10.68%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "1.18"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Value = "0.0888"
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("D45").Value = "90.96"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "1.374.12"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "7.44"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "45.73"
$ws.Range("E51").Value = "  +3.26%  "
